$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is being updated to reflect FY2018 (instead of FY2019) figures.
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 1015849038.23
$ws.Range("P2").Value = 399922409.81
$ws.Range("Q2").Value = 308157654.68
$ws.Range("R2").Value = 32.2192068697
$ws.Range("S2").Value = 40686811.66
$ws.Range("T2").Value = 29.815987763
$ws.Range("U2").Value = 48535906.06
$ws.Range("V2").Value = 26.241060205
$ws.Range("W2").Value = 528359593.55
$ws.Range("X2").Value = 149816868.15
$ws.Range("Y2").Value = 4.9200726948
$ws.Range("Z2").Value = 196070095.18
$ws.Range("AA2").Value = 2.6806481621
$ws.Range("AB2").Value = 487489444.68
$ws.Range("AC2").Value = 22.2065662002
$ws.Range("AD2").Value = 16.713372109
$ws.Range("AE2").Value = 12.0656666135
$ws.Range("AF2").Value = 89.0335584696
$ws.Range("AG2").Value = 52.0116251201
